$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sdflkj"
$ws.Range("B2").Value = "sdfsd"
$ws.Range("A3").Value = "sdflkj3333"
$ws.Range("B3").Value = "sdfsd33"
